# Update commonness index calculations (H and I columns, rows 2-16)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ H = 0.16301143169403076;   I = 0.5299953818321228 }
    3  = @{ H = 0.17079822719097137;   I = 0.52942371368408203 }
    4  = @{ H = 0.20533210039138794;   I = 0.54037958383560181 }
    5  = @{ H = 0.85838252305984497;   I = 0.93988436460494995 }
    6  = @{ H = 0.85841119289398193;   I = 0.9361426830291748 }
    7  = @{ H = 0.85799229145050049;   I = 0.93365299701690674 }
    8  = @{ H = 0.19332195818424225;   I = 0.34843769669532776 }
    9  = @{ H = 0.20416003465652466;   I = 0.34839457273483276 }
    10 = @{ H = 0.20557700097560883;   I = 0.35506138205528259 }
    11 = @{ H = 0.13654132187366486;   I = 0.45524618029594421 }
    12 = @{ H = 0.14872801303863525;   I = 0.45659786462783813 }
    13 = @{ H = 0.18394468724727631;   I = 0.46406301856040955 }
    14 = @{ H = -0.0065314383246004581; I = 0.13369239866733551 }
    15 = @{ H = -0.015694240108132362; I = 0.12246883660554886 }
    16 = @{ H = -0.012039145454764366; I = 0.12394032627344131 }
}

foreach ($row in $values.Keys) {
    $ws.Range("H$row").Value = $values[$row].H
    $ws.Range("I$row").Value = $values[$row].I
}
